$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("F3").Value = ""
$ws.Range("K3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4
$ws.Range("F4").Value = ""
$ws.Range("K4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5
$ws.Range("F5").Value = ""
$ws.Range("K5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6
$ws.Range("K6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 8 - shift values: F becomes empty, K gets old F value
$ws.Range("F8").Value = ""
$ws.Range("K8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9 - shift values: F becomes empty, K gets old F value
$ws.Range("F9").Value = ""
$ws.Range("K9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10 - shift values: F becomes empty, K gets old F value
$ws.Range("F10").Value = ""
$ws.Range("K10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 11
$ws.Range("K11").Value = "3,05 TL - 6,09 TL - 76,17 TL"

# Row 12
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 404,16 TL–3.403,42 TL"

# Row 13
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("F13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14
$ws.Range("F14").Value = ""
$ws.Range("K14").Value = "914,14 TL - 4.265,98 TL"
